$d = $word.ActiveDocument

# Locate the paragraph containing the "{m:userdoc 'zone1'}" field marker and
# select its whole text so we can rebuild it as four separate runs:
#   "{" / "m" / ":userdoc 'zone1'" / "}"
# (matching the TokenIteratorFieldRewriterSplit behaviour that splits a
# field's opening brace, tag name, body and closing brace into distinct runs).
$target = $d.Content
$found = $target.Find.Execute("{m:userdoc 'zone1'}", $false, $false, $false,
                               $false, $false, $true, 1, $false, "", 0)

$target.Delete()

$openDelim = [char]0x7B
$closeDelim = [char]0x7D
$quote = [char]0x27

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
       '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
       '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
       '<pkg:xmlData>' +
       '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
       '<w:body>' +
       '<w:p>' +
       '<w:r><w:t>' + $openDelim + '</w:t></w:r>' +
       '<w:r><w:t>m</w:t></w:r>' +
       '<w:r><w:t>:userdoc ' + $quote + 'zone1' + $quote + '</w:t></w:r>' +
       '<w:r><w:t xml:space="preserve">' + $closeDelim + '</w:t></w:r>' +
       '</w:p>' +
       '</w:body>' +
       '</w:document>' +
       '</pkg:xmlData></pkg:part></pkg:package>'

$target.InsertXML($xml)
